$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current blank row 5, pushing everything
# (the blank spacer row, the header row and all data rows) down by one.
$ws.Rows.Item(5).Insert(-4121, 0)

# Populate the newly inserted row with the outliers footnote, styled to
# match the sheet's existing "Times Roman" body font but left-aligned
# (a brand new combination of already-existing font #2 + left alignment).
$cell = $ws.Range("A5")
$cell.Value2 = "Outliers were removed prior to data analysis. Outliers from SY1140A (2 for Leaf #, 1 for Leaf length, and 2 for Root length)."
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4107
$cell.Font.Name = "Times Roman"
$cell.Font.Family = 2

# The data table (now at A40:N47, previously A39:N46) keeps its sort
# definition; re-apply it so the persisted sortState range follows the
# shift caused by the inserted row.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B40:B47"))
$sortObj.SetRange($ws.Range("A40:N47"))
$sortObj.Header = -4142
$sortObj.Apply()

# Restore the author's last selection.
$ws.Range("E2").Select()
